# New crime data collected - weekly update for week covering 10/16/2023 - 10/22/2023

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: Volume/Number and the reporting week date range
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# ---------------------------------------------------------------------------
# Helper: convert a numeric cell to the "dash / undefined" text placeholder
# cells used throughout the sheet ("0" for counts, "***.*" for % change),
# reusing the formatting of an existing placeholder cell so no new styles
# or shared strings get introduced.
# ---------------------------------------------------------------------------
function Set-PlaceholderText($targetAddr, $templateAddr, $text) {
    $ws.Range($templateAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($targetAddr).Formula = "=""" + $text + """"
    $ws.Range($targetAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# ---------------------------------------------------------------------------
# Row 16 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -75
$ws.Range("I16").Value = 18
$ws.Range("J16").Value = 25
$ws.Range("K16").Value = -28
$ws.Range("L16").Value = 38.461538461538
$ws.Range("M16").Value = -18.181818181818
$ws.Range("N16").Value = -89.473684210526

# ---------------------------------------------------------------------------
# Row 18 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("F18").Value = 1
$ws.Range("N18").Value = -84.375

# ---------------------------------------------------------------------------
# Row 19 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 1
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 25
$ws.Range("J19").Value = 23
$ws.Range("K19").Value = 91.304347826087
$ws.Range("L19").Value = 100
$ws.Range("N19").Value = -71.612903225806

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = -66.666666666666
$ws.Range("F21").Value = 9
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 77
$ws.Range("J21").Value = 66
$ws.Range("K21").Value = 16.666666666666
$ws.Range("L21").Value = 54
$ws.Range("M21").Value = -17.204301075268
$ws.Range("N21").Value = -81.219512195122

# ---------------------------------------------------------------------------
# Row 24 (Shooting Vic.)
# ---------------------------------------------------------------------------
Set-PlaceholderText "C24" "C22" "0"
Set-PlaceholderText "D24" "D22" "0"
Set-PlaceholderText "E24" "E22" "***.*"
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 3
$ws.Range("H24").Value = -33.333333333333
$ws.Range("L24").Value = 17.857142857142
$ws.Range("M24").Value = -59.756097560975

# ---------------------------------------------------------------------------
# Row 25 (Shooting Inc.)
# ---------------------------------------------------------------------------
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("C25").Value = 1
Set-PlaceholderText "D25" "D22" "0"
Set-PlaceholderText "E25" "E22" "***.*"
$ws.Range("F25").Value = 3
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 44
$ws.Range("K25").Value = 62.962962962963
$ws.Range("L25").Value = 109.52380952381
$ws.Range("M25").Value = 158.823529411765

# ---------------------------------------------------------------------------
# Column E width (auto best-fit width shrank once E16 text got shorter)
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 6.7
